# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 249 of the data table, shifting the
# existing rows 249:293 down to 250:294 (same as the canonical OOXML diff,
# which renumbers every row from 249 onward by +1 and appends one new data
# row at the bottom of the shifted range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 249:293 down one row, creating a blank row 249.
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row with the new weekly observation. All the
# "template" columns (market/region/category/unit/etc.) are identical to the
# rest of the block; only the date, volume and the three prices differ.
$ws.Cells.Item(249, 1).Value = 10
$ws.Cells.Item(249, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(249, 3).Value = "La Araucanía"
$ws.Cells.Item(249, 4).Value = 44694
$ws.Cells.Item(249, 5).Value = 9
$ws.Cells.Item(249, 6).Value = 100112017
$ws.Cells.Item(249, 7).Value = "Apio"
$ws.Cells.Item(249, 8).Value = "Americana (o)"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 65
$ws.Cells.Item(249, 11).Value = 10000
$ws.Cells.Item(249, 12).Value = 10000
$ws.Cells.Item(249, 13).Value = 10000
$ws.Cells.Item(249, 14).Value = "`$/docena de matas"
$ws.Cells.Item(249, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(249, 16).Value = 1667
$ws.Cells.Item(249, 17).Value = 6
$ws.Cells.Item(249, 18).Value = "Hortaliza"
